# Scheduled runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, Leve Price*/Profit* columns) across several sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 249
$ws.Range("I6").Value = 249
$ws.Range("K6").Value = 747
$ws.Range("M6").Value = -635
$ws.Range("H8").Value = 3892.5
$ws.Range("I8").Value = 3892.5
$ws.Range("K8").Value = 11677.5
$ws.Range("M8").Value = -11538.5
$ws.Range("H16").Value = 4725
$ws.Range("I16").Value = 3700
$ws.Range("J16").Value = 5750
$ws.Range("K16").Value = 3700
$ws.Range("L16").Value = 5750
$ws.Range("M16").Value = -3470
$ws.Range("N16").Value = -6210
$ws.Range("H28").Value = 5970.9546
$ws.Range("I28").Value = 2647.5386
$ws.Range("K28").Value = 2647.5386
$ws.Range("M28").Value = -2162.5386
$ws.Range("H31").Value = 4082.3333
$ws.Range("I31").Value = 624
$ws.Range("K31").Value = 1872
$ws.Range("M31").Value = -1642
$ws.Range("H39").Value = 740.4737
$ws.Range("I39").Value = 111.545456
$ws.Range("J39").Value = 1605.25
$ws.Range("K39").Value = 334.636368
$ws.Range("L39").Value = 4815.75
$ws.Range("M39").Value = -38.636368
$ws.Range("N39").Value = -5407.75
$ws.Range("H40").Value = 71431020
$ws.Range("I40").Value = 2491
$ws.Range("J40").Value = 125002420
$ws.Range("K40").Value = 2491
$ws.Range("L40").Value = 125002420
$ws.Range("M40").Value = -2316
$ws.Range("N40").Value = -125002770
$ws.Range("H43").Value = 5500
$ws.Range("I43").Value = 5500
$ws.Range("K43").Value = 5500
$ws.Range("M43").Value = -5431
$ws.Range("H62").Value = 6086.375
$ws.Range("I62").Value = 2098.2
$ws.Range("K62").Value = 2098.2
$ws.Range("M62").Value = -1474.2
$ws.Range("H65").Value = 6086.375
$ws.Range("I65").Value = 2098.2
$ws.Range("K65").Value = 10491
$ws.Range("M65").Value = -7371
$ws.Range("H100").Value = 8631.666999999999
$ws.Range("I100").Value = 6345
$ws.Range("J100").Value = 10461
$ws.Range("K100").Value = 6345
$ws.Range("L100").Value = 10461
$ws.Range("M100").Value = -5804
$ws.Range("N100").Value = -11543
$ws.Range("H132").Value = 7972.636
$ws.Range("I132").Value = 4001
$ws.Range("K132").Value = 12003
$ws.Range("M132").Value = -9473
$ws.Range("H135").Value = 1760.6757
$ws.Range("I135").Value = 488.53333
$ws.Range("K135").Value = 4396.79997
$ws.Range("M135").Value = -1861.79997
$ws.Range("H137").Value = 928526.5
$ws.Range("I137").Value = 799.4
$ws.Range("J137").Value = 3247844.2
$ws.Range("K137").Value = 2398.2
$ws.Range("L137").Value = 9743532.600000001
$ws.Range("M137").Value = 151.8000000000002
$ws.Range("N137").Value = -9748632.600000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2598.8462
$ws.Range("I32").Value = 2554.4644
$ws.Range("J32").Value = 2875
$ws.Range("K32").Value = 2554.4644
$ws.Range("L32").Value = 2875
$ws.Range("M32").Value = -2267.4644
$ws.Range("N32").Value = -3449
$ws.Range("H45").Value = 3187
$ws.Range("I45").Value = 1450.625
$ws.Range("J45").Value = 5965.2
$ws.Range("K45").Value = 1450.625
$ws.Range("L45").Value = 5965.2
$ws.Range("M45").Value = -1073.625
$ws.Range("N45").Value = -6719.2
$ws.Range("H60").Value = 66644.42999999999
$ws.Range("I60").Value = 66644.42999999999
$ws.Range("K60").Value = 66644.42999999999
$ws.Range("M60").Value = -65911.42999999999
$ws.Range("H110").Value = 3306.1
$ws.Range("I110").Value = 2562.3333
$ws.Range("K110").Value = 2562.3333
$ws.Range("M110").Value = -517.3332999999998
$ws.Range("H122").Value = 3578.238
$ws.Range("I122").Value = 2407
$ws.Range("K122").Value = 7221
$ws.Range("M122").Value = -4771
$ws.Range("H132").Value = 5269243.5
$ws.Range("I132").Value = 6151.7856
$ws.Range("K132").Value = 18455.3568
$ws.Range("M132").Value = -15925.3568
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 14294147
$ws.Range("I16").Value = 33340004
$ws.Range("K16").Value = 33340004
$ws.Range("M16").Value = -33339717
$ws.Range("H113").Value = 14294147
$ws.Range("I113").Value = 33340004
$ws.Range("K113").Value = 33340004
$ws.Range("M113").Value = -33337834
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1168.5
$ws.Range("I113").Value = 853.2308
$ws.Range("J113").Value = 1623.8889
$ws.Range("K113").Value = 2559.6924
$ws.Range("L113").Value = 4871.6667
$ws.Range("M113").Value = -389.6923999999999
$ws.Range("N113").Value = -9211.6667
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3156
$ws.Range("I102").Value = 3156
$ws.Range("K102").Value = 3156
$ws.Range("M102").Value = -1534
$ws.Range("H122").Value = 3968.8
$ws.Range("I122").Value = 4064.1538
$ws.Range("K122").Value = 12192.4614
$ws.Range("M122").Value = -9742.4614
$ws.Range("H126").Value = 2735.6
$ws.Range("J126").Value = 2419.5
$ws.Range("L126").Value = 7258.5
$ws.Range("N126").Value = -12198.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 299.66666
$ws.Range("I22").Value = 199
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 199
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = 96
$ws.Range("N22").Value = -940
$ws.Range("H27").Value = 299.66666
$ws.Range("I27").Value = 199
$ws.Range("J27").Value = 350
$ws.Range("K27").Value = 199
$ws.Range("L27").Value = 350
$ws.Range("M27").Value = -92
$ws.Range("N27").Value = -564
$ws.Range("H46").Value = 1591.6666
$ws.Range("J46").Value = 1390
$ws.Range("L46").Value = 1390
$ws.Range("N46").Value = -1766
$ws.Range("H68").Value = 5210852
$ws.Range("I68").Value = 10418531
$ws.Range("J68").Value = 3173.25
$ws.Range("K68").Value = 10418531
$ws.Range("L68").Value = 3173.25
$ws.Range("M68").Value = -10417782
$ws.Range("N68").Value = -4671.25
$ws.Range("H71").Value = 5210852
$ws.Range("I71").Value = 10418531
$ws.Range("J71").Value = 3173.25
$ws.Range("K71").Value = 52092655
$ws.Range("L71").Value = 15866.25
$ws.Range("M71").Value = -52088911
$ws.Range("N71").Value = -23354.25
$ws.Range("H100").Value = 31286106
$ws.Range("I100").Value = 4974.75
$ws.Range("J100").Value = 62567236
$ws.Range("K100").Value = 4974.75
$ws.Range("L100").Value = 62567236
$ws.Range("M100").Value = -4433.75
$ws.Range("N100").Value = -62568318
$ws.Range("H132").Value = 3894.9546
$ws.Range("I132").Value = 2606.0625
$ws.Range("K132").Value = 7818.1875
$ws.Range("M132").Value = -5288.1875
$ws.Range("H136").Value = 3613.2856
$ws.Range("I136").Value = 2679.625
$ws.Range("K136").Value = 8038.875
$ws.Range("M136").Value = -5488.875
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1629.6
$ws.Range("I100").Value = 1599.5
$ws.Range("K100").Value = 3199
$ws.Range("M100").Value = -2658
$ws.Range("H126").Value = 7313.8
$ws.Range("I126").Value = 8301.08
$ws.Range("J126").Value = 2377.4
$ws.Range("K126").Value = 24903.24
$ws.Range("L126").Value = 7132.200000000001
$ws.Range("M126").Value = -22433.24
$ws.Range("N126").Value = -12072.2
$ws.Range("H132").Value = 528917.9399999999
$ws.Range("I132").Value = 2818
$ws.Range("K132").Value = 8454
$ws.Range("M132").Value = -5924
